$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before GM, shifting GM->GN (nom) and GN->GO (url_produit) to the right.
$ws.Columns("GM").Insert()

# New timestamp header for the freshly inserted column.
$ws.Range("GM1").Value = "2026-02-05 23:14:58"

# Rows 2-80 previously had their last recorded price in column GL (now unchanged after the insert,
# since the insert only shifted GM onward). Carry that same price forward into the new GM column
# (matches the "duplicate the prior price under the new timestamp" pattern in the diff).
$ws.Range("GM2").Value = 39.83
$ws.Range("GM3").Value = 169.95
$ws.Range("GM4").Value = 249.95
$ws.Range("GM5").Value = 299.95
$ws.Range("GM6").Value = 339.95
$ws.Range("GM7").Value = 619
$ws.Range("GM8").Value = 619
$ws.Range("GM9").Value = 659
$ws.Range("GM10").Value = 659
$ws.Range("GM11").Value = 749
$ws.Range("GM12").Value = 809
$ws.Range("GM13").Value = 809
$ws.Range("GM14").Value = 809
$ws.Range("GM15").Value = 809
$ws.Range("GM16").Value = 809
$ws.Range("GM17").Value = 849
$ws.Range("GM18").Value = 899
$ws.Range("GM19").Value = 899
$ws.Range("GM20").Value = 909
$ws.Range("GM21").Value = 909
$ws.Range("GM22").Value = 909
$ws.Range("GM23").Value = 909
$ws.Range("GM24").Value = 909
$ws.Range("GM25").Value = 969
$ws.Range("GM26").Value = 969
$ws.Range("GM27").Value = 969
$ws.Range("GM28").Value = 969
$ws.Range("GM29").Value = 969
$ws.Range("GM30").Value = 999
$ws.Range("GM31").Value = 999
$ws.Range("GM32").Value = 1039
$ws.Range("GM33").Value = 1039
$ws.Range("GM34").Value = 1079
$ws.Range("GM35").Value = 1079
$ws.Range("GM36").Value = 1079
$ws.Range("GM37").Value = 1079
$ws.Range("GM38").Value = 1099
$ws.Range("GM39").Value = 1099
$ws.Range("GM40").Value = 1199
$ws.Range("GM41").Value = 1219
$ws.Range("GM42").Value = 1219
$ws.Range("GM43").Value = 1219
$ws.Range("GM44").Value = 1219
$ws.Range("GM45").Value = 1219
$ws.Range("GM46").Value = 1229
$ws.Range("GM47").Value = 1229
$ws.Range("GM48").Value = 1249
$ws.Range("GM49").Value = 1329
$ws.Range("GM50").Value = 1329
$ws.Range("GM51").Value = 1329
$ws.Range("GM52").Value = 1329
$ws.Range("GM53").Value = 1329
$ws.Range("GM54").Value = 1329
$ws.Range("GM55").Value = 1329
$ws.Range("GM56").Value = 1349
$ws.Range("GM57").Value = 1419
$ws.Range("GM58").Value = 1479
$ws.Range("GM59").Value = 1479
$ws.Range("GM60").Value = 1479
$ws.Range("GM61").Value = 1549
$ws.Range("GM62").Value = 1579
$ws.Range("GM63").Value = 1579
$ws.Range("GM64").Value = 1579
$ws.Range("GM65").Value = 1579
$ws.Range("GM66").Value = 1579
$ws.Range("GM67").Value = 1579
$ws.Range("GM68").Value = 1579
$ws.Range("GM69").Value = 1729
$ws.Range("GM70").Value = 1729
$ws.Range("GM71").Value = 1729
$ws.Range("GM72").Value = 1829
$ws.Range("GM73").Value = 1829
$ws.Range("GM74").Value = 1829
$ws.Range("GM75").Value = 1979
$ws.Range("GM76").Value = 1979
$ws.Range("GM77").Value = 1979
$ws.Range("GM78").Value = 2479
$ws.Range("GM79").Value = 2479
$ws.Range("GM80").Value = 2479
